$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Record progress update: Anushka Chincholkar (row 4) now has a "Work Done Upto" entry.
$ws.Range("B4").Value = "Studied basic concepts of development in flutter"

# Excel leaves the edited cell as the active selection after data entry.
$ws.Range("B4").Select()
